# Auto-generated edit script applying Anima_Profits.xlsx leve-profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 141
$ws.Cells.Item(141, 8).Value = 7163.222  # H141
$ws.Cells.Item(141, 9).Value = 3397.3333  # I141
$ws.Cells.Item(141, 10).Value = 14695  # J141
$ws.Cells.Item(141, 11).Value = 10191.9999  # K141
$ws.Cells.Item(141, 12).Value = 44085  # L141
$ws.Cells.Item(141, 13).Value = -5011.999899999999  # M141
$ws.Cells.Item(141, 14).Value = -54445  # N141

$ws = $wb.Worksheets.Item("ARM")
# ARM row 12
$ws.Cells.Item(12, 8).Value = 600  # H12
$ws.Cells.Item(12, 9).Value = 600  # I12
$ws.Cells.Item(12, 11).Value = 600  # K12
$ws.Cells.Item(12, 13).Value = -427  # M12

# ARM row 32
$ws.Cells.Item(32, 8).Value = 1492753.1  # H32
$ws.Cells.Item(32, 9).Value = 1860490.8  # I32
$ws.Cells.Item(32, 11).Value = 1860490.8  # K32
$ws.Cells.Item(32, 13).Value = -1860203.8  # M32

# ARM row 45
$ws.Cells.Item(45, 8).Value = 1667.52  # H45
$ws.Cells.Item(45, 9).Value = 1374.7646  # I45
$ws.Cells.Item(45, 10).Value = 2289.625  # J45
$ws.Cells.Item(45, 11).Value = 1374.7646  # K45
$ws.Cells.Item(45, 12).Value = 2289.625  # L45
$ws.Cells.Item(45, 13).Value = -997.7646  # M45
$ws.Cells.Item(45, 14).Value = -3043.625  # N45

# ARM row 81
$ws.Cells.Item(81, 8).Value = 0  # H81
$ws.Cells.Item(81, 10).Value = 0  # J81
$ws.Cells.Item(81, 12).Value = 0  # L81
$ws.Cells.Item(81, 14).ClearContents()  # N81

# ARM row 84
$ws.Cells.Item(84, 8).Value = 0  # H84
$ws.Cells.Item(84, 10).Value = 0  # J84
$ws.Cells.Item(84, 12).Value = 0  # L84
$ws.Cells.Item(84, 14).ClearContents()  # N84

# ARM row 102
$ws.Cells.Item(102, 8).Value = 3000  # H102
$ws.Cells.Item(102, 9).Value = 3000  # I102
$ws.Cells.Item(102, 10).Value = 0  # J102
$ws.Cells.Item(102, 11).Value = 3000  # K102
$ws.Cells.Item(102, 12).Value = 0  # L102
$ws.Cells.Item(102, 13).Value = -1378  # M102
$ws.Cells.Item(102, 14).ClearContents()  # N102

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Cells.Item(31, 8).Value = 5241.3667  # H31
$ws.Cells.Item(31, 9).Value = 1405.5264  # I31
$ws.Cells.Item(31, 10).Value = 7018.951  # J31
$ws.Cells.Item(31, 11).Value = 1405.5264  # K31
$ws.Cells.Item(31, 12).Value = 7018.951  # L31
$ws.Cells.Item(31, 13).Value = -1110.5264  # M31
$ws.Cells.Item(31, 14).Value = -7608.951  # N31

# CRP row 34
$ws.Cells.Item(34, 8).Value = 5241.3667  # H34
$ws.Cells.Item(34, 9).Value = 1405.5264  # I34
$ws.Cells.Item(34, 10).Value = 7018.951  # J34
$ws.Cells.Item(34, 11).Value = 1405.5264  # K34
$ws.Cells.Item(34, 12).Value = 7018.951  # L34
$ws.Cells.Item(34, 13).Value = -1203.5264  # M34
$ws.Cells.Item(34, 14).Value = -7422.951  # N34

# CRP row 81
$ws.Cells.Item(81, 8).Value = 93664  # H81
$ws.Cells.Item(81, 10).Value = 93664  # J81
$ws.Cells.Item(81, 12).Value = 93664  # L81
$ws.Cells.Item(81, 14).Value = -95660  # N81

# CRP row 84
$ws.Cells.Item(84, 8).Value = 93664  # H84
$ws.Cells.Item(84, 10).Value = 93664  # J84
$ws.Cells.Item(84, 12).Value = 280992  # L84
$ws.Cells.Item(84, 14).Value = -290976  # N84

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Cells.Item(5, 8).Value = 707.65515  # H5
$ws.Cells.Item(5, 9).Value = 503.91666  # I5
$ws.Cells.Item(5, 10).Value = 1041.0454  # J5
$ws.Cells.Item(5, 11).Value = 1511.74998  # K5
$ws.Cells.Item(5, 12).Value = 3123.1362  # L5
$ws.Cells.Item(5, 13).Value = -1399.74998  # M5
$ws.Cells.Item(5, 14).Value = -3347.1362  # N5

# CUL row 62
$ws.Cells.Item(62, 8).Value = 7838.8335  # H62
$ws.Cells.Item(62, 9).Value = 1012  # I62
$ws.Cells.Item(62, 10).Value = 8459.454  # J62
$ws.Cells.Item(62, 11).Value = 3036  # K62
$ws.Cells.Item(62, 12).Value = 25378.362  # L62
$ws.Cells.Item(62, 13).Value = -2350  # M62
$ws.Cells.Item(62, 14).Value = -26750.362  # N62

# CUL row 65
$ws.Cells.Item(65, 8).Value = 7838.8335  # H65
$ws.Cells.Item(65, 9).Value = 1012  # I65
$ws.Cells.Item(65, 10).Value = 8459.454  # J65
$ws.Cells.Item(65, 11).Value = 9108  # K65
$ws.Cells.Item(65, 12).Value = 76135.086  # L65
$ws.Cells.Item(65, 13).Value = -5676  # M65
$ws.Cells.Item(65, 14).Value = -82999.086  # N65

# CUL row 68
$ws.Cells.Item(68, 8).Value = 1266.1428  # H68
$ws.Cells.Item(68, 9).Value = 1375.5  # I68
$ws.Cells.Item(68, 10).Value = 1222.4  # J68
$ws.Cells.Item(68, 11).Value = 4126.5  # K68
$ws.Cells.Item(68, 12).Value = 3667.2  # L68
$ws.Cells.Item(68, 13).Value = -3315.5  # M68
$ws.Cells.Item(68, 14).Value = -5289.200000000001  # N68

# CUL row 71
$ws.Cells.Item(71, 8).Value = 1266.1428  # H71
$ws.Cells.Item(71, 9).Value = 1375.5  # I71
$ws.Cells.Item(71, 10).Value = 1222.4  # J71
$ws.Cells.Item(71, 11).Value = 12379.5  # K71
$ws.Cells.Item(71, 12).Value = 11001.6  # L71
$ws.Cells.Item(71, 13).Value = -8323.5  # M71
$ws.Cells.Item(71, 14).Value = -19113.6  # N71

# CUL row 74
$ws.Cells.Item(74, 8).Value = 1756.5  # H74
$ws.Cells.Item(74, 9).Value = 1008.6667  # I74
$ws.Cells.Item(74, 10).Value = 4000  # J74
$ws.Cells.Item(74, 11).Value = 3026.0001  # K74
$ws.Cells.Item(74, 12).Value = 12000  # L74
$ws.Cells.Item(74, 13).Value = -1965.0001  # M74
$ws.Cells.Item(74, 14).Value = -14122  # N74

# CUL row 75
$ws.Cells.Item(75, 8).Value = 999  # H75
$ws.Cells.Item(75, 9).Value = 999  # I75
$ws.Cells.Item(75, 11).Value = 2997  # K75
$ws.Cells.Item(75, 13).Value = -1999  # M75

# CUL row 77
$ws.Cells.Item(77, 8).Value = 1756.5  # H77
$ws.Cells.Item(77, 9).Value = 1008.6667  # I77
$ws.Cells.Item(77, 10).Value = 4000  # J77
$ws.Cells.Item(77, 11).Value = 9078.0003  # K77
$ws.Cells.Item(77, 12).Value = 36000  # L77
$ws.Cells.Item(77, 13).Value = -3774.0003  # M77
$ws.Cells.Item(77, 14).Value = -46608  # N77

# CUL row 78
$ws.Cells.Item(78, 8).Value = 999  # H78
$ws.Cells.Item(78, 9).Value = 999  # I78
$ws.Cells.Item(78, 11).Value = 8991  # K78
$ws.Cells.Item(78, 13).Value = -3999  # M78

# CUL row 82
$ws.Cells.Item(82, 8).Value = 1000  # H82
$ws.Cells.Item(82, 9).Value = 1000  # I82
$ws.Cells.Item(82, 11).Value = 3000  # K82
$ws.Cells.Item(82, 13).Value = -2594  # M82

# CUL row 85
$ws.Cells.Item(85, 8).Value = 1000  # H85
$ws.Cells.Item(85, 9).Value = 1000  # I85
$ws.Cells.Item(85, 11).Value = 3000  # K85
$ws.Cells.Item(85, 13).Value = -1596  # M85

# CUL row 135
$ws.Cells.Item(135, 8).Value = 707.65515  # H135
$ws.Cells.Item(135, 9).Value = 503.91666  # I135
$ws.Cells.Item(135, 10).Value = 1041.0454  # J135
$ws.Cells.Item(135, 11).Value = 4535.24994  # K135
$ws.Cells.Item(135, 12).Value = 9369.408599999999  # L135
$ws.Cells.Item(135, 13).Value = -2000.24994  # M135
$ws.Cells.Item(135, 14).Value = -14439.4086  # N135

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70
$ws.Cells.Item(70, 8).Value = 9595.817999999999  # H70
$ws.Cells.Item(70, 9).Value = 10723.529  # I70
$ws.Cells.Item(70, 10).Value = 5761.6  # J70
$ws.Cells.Item(70, 11).Value = 10723.529  # K70
$ws.Cells.Item(70, 12).Value = 5761.6  # L70
$ws.Cells.Item(70, 13).Value = -10453.529  # M70
$ws.Cells.Item(70, 14).Value = -6301.6  # N70

# GSM row 73
$ws.Cells.Item(73, 8).Value = 9595.817999999999  # H73
$ws.Cells.Item(73, 9).Value = 10723.529  # I73
$ws.Cells.Item(73, 10).Value = 5761.6  # J73
$ws.Cells.Item(73, 11).Value = 10723.529  # K73
$ws.Cells.Item(73, 12).Value = 5761.6  # L73
$ws.Cells.Item(73, 13).Value = -9787.529  # M73
$ws.Cells.Item(73, 14).Value = -7633.6  # N73

# GSM row 126
$ws.Cells.Item(126, 8).Value = 3240.9333  # H126
$ws.Cells.Item(126, 9).Value = 2700.25  # I126
$ws.Cells.Item(126, 11).Value = 8100.75  # K126
$ws.Cells.Item(126, 13).Value = -5630.75  # M126

# GSM row 132
$ws.Cells.Item(132, 8).Value = 3432.7727  # H132
$ws.Cells.Item(132, 9).Value = 2565.6155  # I132
$ws.Cells.Item(132, 10).Value = 4685.3335  # J132
$ws.Cells.Item(132, 11).Value = 7696.8465  # K132
$ws.Cells.Item(132, 12).Value = 14056.0005  # L132
$ws.Cells.Item(132, 13).Value = -5166.8465  # M132
$ws.Cells.Item(132, 14).Value = -19116.0005  # N132

# GSM row 136
$ws.Cells.Item(136, 8).Value = 9108.174000000001  # H136
$ws.Cells.Item(136, 9).Value = 0  # I136
$ws.Cells.Item(136, 10).Value = 9108.174000000001  # J136
$ws.Cells.Item(136, 11).Value = 0  # K136
$ws.Cells.Item(136, 12).Value = 27324.522  # L136
$ws.Cells.Item(136, 13).ClearContents()  # M136
$ws.Cells.Item(136, 14).Value = -32424.522  # N136

$ws = $wb.Worksheets.Item("LTW")
# LTW row 56
$ws.Cells.Item(56, 8).Value = 19900  # H56
$ws.Cells.Item(56, 9).Value = 19900  # I56
$ws.Cells.Item(56, 11).Value = 19900  # K56
$ws.Cells.Item(56, 13).Value = -19209  # M56

# LTW row 115
$ws.Cells.Item(115, 8).Value = 40555.5  # H115
$ws.Cells.Item(115, 10).Value = 40555.5  # J115
$ws.Cells.Item(115, 12).Value = 40555.5  # L115
$ws.Cells.Item(115, 14).Value = -42905.5  # N115

$ws = $wb.Worksheets.Item("WVR")
# WVR row 110
$ws.Cells.Item(110, 8).Value = 55644  # H110
$ws.Cells.Item(110, 10).Value = 55644  # J110
$ws.Cells.Item(110, 12).Value = 55644  # L110
$ws.Cells.Item(110, 14).Value = -63824  # N110

# WVR row 116
$ws.Cells.Item(116, 8).Value = 0  # H116
$ws.Cells.Item(116, 10).Value = 0  # J116
$ws.Cells.Item(116, 12).Value = 0  # L116
$ws.Cells.Item(116, 14).ClearContents()  # N116
